$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 7).Value = 1
$ws.Cells.Item(10, 8).Value = 2

$ws.Cells.Item(15, 5).Value = 157
$ws.Cells.Item(15, 7).Value = 41
$ws.Cells.Item(15, 8).Value = 123

$ws.Cells.Item(26, 5).Value = 27

$ws.Cells.Item(29, 5).Value = 17

$ws.Cells.Item(68, 5).Value = 17

$ws.Cells.Item(73, 5).Value = 28
$ws.Cells.Item(73, 6).Value = 10
$ws.Cells.Item(73, 8).Value = 22

$ws.Cells.Item(74, 5).Value = 19

$ws.Cells.Item(77, 5).Value = 53
$ws.Cells.Item(77, 6).Value = 20
$ws.Cells.Item(77, 8).Value = 37
